# Correção de erro de digitação:
# O campo "Perfil" (F20) na aba "FIIs - Análise Investimeno" estava definido
# como "Agressivo" e deveria ser "Moderado". Também restauramos a aba
# "FIIs - Análise Investimeno" como aba ativa/selecionada, com a visão
# rolada de volta para o topo (célula C1).

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("FIIs - Análise Investimeno")
$wsPerfis = $wb.Worksheets.Item("Perfis")

# Corrige o valor selecionado do perfil de investimento
$wsMain.Range("F20").Value = "Moderado"

# Recalcula a pasta de trabalho para atualizar formulas dependentes (E57:E62,
# F57:F62, etc.) e os caches dos gráficos.
$excel.CalculateFullRebuild()

# Ativa a planilha principal e garante que a seleção/rolagem fique no topo.
$wsMain.Activate()
$wsMain.Range("F20").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
